$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of test case data (row 3), mirroring row 2's structure
$ws.Range("A3").Value = "TC_AR_002"
$ws.Range("B3").Value = "New Account Registration"
$ws.Range("C3").Value = "Lady2"
$ws.Range("D3").Value = "Raga2"
$ws.Range("E3").Value = "raga2@yopmail.com"
$ws.Range("F3").Value = "Smart@1234"
$ws.Range("G3").Value = "Smart@1234"

# Apply same style as row 2's hyperlink cells to row 3
$ws.Range("E3:G3").Style = "Hyperlink"

# Add hyperlinks mirroring the ones on row 2
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:raga2@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:Smart@1234")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:Smart@1234")

# Update selection to mirror the edited cell
$ws.Range("C4").Select()
